$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 60145.65
$ws.Range("J28").Value = 318.6
$ws.Range("L28").Value = 318.6
$ws.Range("N28").Value = -1288.6
$ws.Range("H80").Value = 5714.952
$ws.Range("I80").Value = 10195.182
$ws.Range("J80").Value = 786.7
$ws.Range("K80").Value = 30585.546
$ws.Range("L80").Value = 2360.1
$ws.Range("M80").Value = -29587.546
$ws.Range("N80").Value = -4356.1
$ws.Range("H83").Value = 5714.952
$ws.Range("I83").Value = 10195.182
$ws.Range("J83").Value = 786.7
$ws.Range("K83").Value = 91756.63800000001
$ws.Range("L83").Value = 7080.3
$ws.Range("M83").Value = -86764.63800000001
$ws.Range("N83").Value = -17064.3
$ws.Range("H106").Value = 2726.6428
$ws.Range("I106").Value = 2726.6428
$ws.Range("K106").Value = 2726.6428
$ws.Range("M106").Value = -2095.6428
$ws.Range("H132").Value = 1719.7858
$ws.Range("I132").Value = 1768.775
$ws.Range("K132").Value = 5306.325000000001
$ws.Range("M132").Value = -2776.325000000001
$ws.Range("H137").Value = 6787.2
$ws.Range("I137").Value = 7484
$ws.Range("K137").Value = 22452
$ws.Range("M137").Value = -19902
$ws.Range("H138").Value = 3972.1755
$ws.Range("J138").Value = 4719.905
$ws.Range("L138").Value = 14159.715
$ws.Range("N138").Value = -24439.715
$ws.Range("H141").Value = 2480.0588
$ws.Range("I141").Value = 2480.0588
$ws.Range("K141").Value = 7440.176399999999
$ws.Range("M141").Value = -2260.176399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1488
$ws.Range("I74").Value = 1515.8823
$ws.Range("K74").Value = 1515.8823
$ws.Range("M74").Value = -641.8823
$ws.Range("H77").Value = 1488
$ws.Range("I77").Value = 1515.8823
$ws.Range("K77").Value = 7579.4115
$ws.Range("M77").Value = -3211.4115
$ws.Range("H122").Value = 4132.2827
$ws.Range("I122").Value = 3297.6333
$ws.Range("K122").Value = 9892.8999
$ws.Range("M122").Value = -7442.8999
$ws.Range("H132").Value = 1850.0422
$ws.Range("I132").Value = 1886.6307
$ws.Range("K132").Value = 5659.8921
$ws.Range("M132").Value = -3129.8921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 296.1
$ws.Range("I22").Value = 282.75
$ws.Range("J22").Value = 349.5
$ws.Range("K22").Value = 282.75
$ws.Range("L22").Value = 349.5
$ws.Range("M22").Value = -109.75
$ws.Range("N22").Value = -695.5
$ws.Range("H94").Value = 750
$ws.Range("J94").Value = 1000
$ws.Range("L94").Value = 1000
$ws.Range("N94").Value = -1902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9137.25
$ws.Range("J16").Value = 17650.2
$ws.Range("L16").Value = 17650.2
$ws.Range("N16").Value = -18224.2
$ws.Range("H113").Value = 9137.25
$ws.Range("J113").Value = 17650.2
$ws.Range("L113").Value = 17650.2
$ws.Range("N113").Value = -21990.2
$ws.Range("H124").Value = 37997.5
$ws.Range("J124").Value = 37997.5
$ws.Range("L124").Value = 37997.5
$ws.Range("N124").Value = -42907.5
$ws.Range("H132").Value = 1674.421
$ws.Range("I132").Value = 1526.25
$ws.Range("J132").Value = 1928.4286
$ws.Range("K132").Value = 4578.75
$ws.Range("L132").Value = 5785.2858
$ws.Range("M132").Value = -2048.75
$ws.Range("N132").Value = -10845.2858
$ws.Range("H134").Value = 335862.28
$ws.Range("I134").Value = 2677.4285
$ws.Range("J134").Value = 5000450
$ws.Range("K134").Value = 8032.2855
$ws.Range("L134").Value = 15001350
$ws.Range("M134").Value = -5497.2855
$ws.Range("N134").Value = -15006420

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 2166941.5
$ws.Range("I9").Value = 550
$ws.Range("K9").Value = 1650
$ws.Range("M9").Value = -1426
$ws.Range("H22").Value = 5598.4
$ws.Range("J22").Value = 5598.4
$ws.Range("L22").Value = 16795.2
$ws.Range("N22").Value = -17133.2
$ws.Range("H27").Value = 5598.4
$ws.Range("J27").Value = 5598.4
$ws.Range("L27").Value = 16795.2
$ws.Range("N27").Value = -16999.2
$ws.Range("H41").Value = 1750
$ws.Range("I41").Value = 1500
$ws.Range("J41").Value = 2000
$ws.Range("K41").Value = 4500
$ws.Range("L41").Value = 6000
$ws.Range("M41").Value = -4162
$ws.Range("N41").Value = -6676
$ws.Range("H44").Value = 1500
$ws.Range("J44").Value = 1500
$ws.Range("L44").Value = 4500
$ws.Range("N44").Value = -5296
$ws.Range("H131").Value = 3916.7693
$ws.Range("J131").Value = 4385
$ws.Range("L131").Value = 13155
$ws.Range("N131").Value = -23235
$ws.Range("H132").Value = 650051.8
$ws.Range("I132").Value = 100867.3
$ws.Range("J132").Value = 1434601.1
$ws.Range("K132").Value = 907805.7000000001
$ws.Range("L132").Value = 12911409.9
$ws.Range("M132").Value = -905275.7000000001
$ws.Range("N132").Value = -12916469.9
$ws.Range("H139").Value = 7585.12
$ws.Range("I139").Value = 7204.846
$ws.Range("K139").Value = 21614.538
$ws.Range("M139").Value = -16474.538

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 9009
$ws.Range("J20").Value = 9009
$ws.Range("L20").Value = 9009
$ws.Range("N20").Value = -9499
$ws.Range("H97").Value = 4577.5
$ws.Range("I97").Value = 3101.5
$ws.Range("J97").Value = 6053.5
$ws.Range("K97").Value = 3101.5
$ws.Range("L97").Value = 6053.5
$ws.Range("M97").Value = -2605.5
$ws.Range("N97").Value = -7045.5
$ws.Range("H125").Value = 91000
$ws.Range("J125").Value = 91000
$ws.Range("L125").Value = 91000
$ws.Range("N125").Value = -95920
$ws.Range("H132").Value = 30559.135
$ws.Range("I132").Value = 4199.0386
$ws.Range("J132").Value = 92864.82000000001
$ws.Range("K132").Value = 12597.1158
$ws.Range("L132").Value = 278594.46
$ws.Range("M132").Value = -10067.1158
$ws.Range("N132").Value = -283654.46
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120
$ws.Range("H137").Value = 49999.832
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 636250
$ws.Range("J20").Value = 636250
$ws.Range("L20").Value = 636250
$ws.Range("N20").Value = -636702
$ws.Range("H93").Value = 2607.9333
$ws.Range("I93").Value = 2047.6666
$ws.Range("K93").Value = 2047.6666
$ws.Range("M93").Value = -799.6666
$ws.Range("H94").Value = 60000
$ws.Range("J94").Value = 60000
$ws.Range("L94").Value = 60000
$ws.Range("N94").Value = -61352
$ws.Range("H122").Value = 910989.75
$ws.Range("I122").Value = 501963.75
$ws.Range("K122").Value = 1505891.25
$ws.Range("M122").Value = -1503441.25
$ws.Range("H136").Value = 718390.5
$ws.Range("I136").Value = 803877.7
$ws.Range("J136").Value = 5997.3335
$ws.Range("K136").Value = 2411633.1
$ws.Range("L136").Value = 17992.0005
$ws.Range("M136").Value = -2409083.1
$ws.Range("N136").Value = -23092.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H81").Value = 1610.1538
$ws.Range("I81").Value = 1175.6364
$ws.Range("K81").Value = 2351.2728
$ws.Range("M81").Value = -1290.2728
$ws.Range("H84").Value = 1610.1538
$ws.Range("I84").Value = 1175.6364
$ws.Range("K84").Value = 11756.364
$ws.Range("M84").Value = -6452.364000000001
$ws.Range("H100").Value = 799.7778
$ws.Range("I100").Value = 774.75
$ws.Range("K100").Value = 1549.5
$ws.Range("M100").Value = -1008.5
$ws.Range("H122").Value = 71432660
$ws.Range("I122").Value = 111114390
$ws.Range("J122").Value = 5541
$ws.Range("K122").Value = 333343170
$ws.Range("L122").Value = 16623
$ws.Range("M122").Value = -333340720
$ws.Range("N122").Value = -21523
$ws.Range("H132").Value = 39749.855
$ws.Range("I132").Value = 2890.389
$ws.Range("K132").Value = 8671.167000000001
$ws.Range("M132").Value = -6141.167000000001
$ws.Range("H136").Value = 8856523
$ws.Range("I136").Value = 10099645
$ws.Range("J136").Value = 403294.6
$ws.Range("K136").Value = 30298935
$ws.Range("L136").Value = 1209883.8
$ws.Range("M136").Value = -30296385
$ws.Range("N136").Value = -1214983.8

